# "Actualización 10 de Mayo" - append rescue-exam ("Rescatables") roster rows.
$ws = $excel.Worksheets.Item("Rescatables")

$A = @(20330051920075, 20330051920061, 20330051920062, 19330051920058, 19330051920058, 19330051920059, 19330051920059, 20330051920088, 20330051920058, 20330051920064, 20330051920066)
$B = @("CONTRERAS", "ROBLES", "ROMERO", "HERNANDEZ", "HERNANDEZ", "HERNANDEZ", "HERNANDEZ", "MAYAHUA", "RAMIREZ", "ROJAS", "SILVESTRE")
$C = @("GARCIA", "IXMATLAHUA", "REYES", "CHAVEZ", "CHAVEZ", "HERNANDEZ", "HERNANDEZ", "TEMOXTLE", "BELLO", "VELASCO", "ARIAS")
$D = @("JORGE HUMBERTO", "ALAN URIEL", "AMANDA MICHEL", "ALEXIS ARMANDO", "ALEXIS ARMANDO", "AGUSTIN", "AGUSTIN", "ELIAS", "ZURIEL ARTURO", "JORGE ALEJANDRO", "YAIR")
$E = @("DISEÑA INSTALACIONES ELÉCTRICAS", "DISEÑA INSTALACIONES ELÉCTRICAS", "DISEÑA INSTALACIONES ELÉCTRICAS", "PROGRAMA Y CONECTA CONTROLADORES LÓGICOS PROGRAMABLES (PLC´S)", "MANTIENE EN OPERACIÓN CIRCUITOS DE CONTROL ELECTROMAGNÉTICO", "PROGRAMA Y CONECTA CONTROLADORES LÓGICOS PROGRAMABLES (PLC´S)", "MANTIENE EN OPERACIÓN CIRCUITOS DE CONTROL ELECTROMAGNÉTICO", "DISEÑA INSTALACIONES ELÉCTRICAS", "DISEÑA INSTALACIONES ELÉCTRICAS", "DISEÑA INSTALACIONES ELÉCTRICAS", "DISEÑA INSTALACIONES ELÉCTRICAS")
$F = @("2AEV", "2AEV", "2AEV", "4AEV", "4AEV", "4AEV", "4AEV", "2AEV", "2AEV", "2AEV", "2AEV")
$G = @(2, 2, 2, 2, 2, 2, 2, 1, 1, 1, 1)

$n = $A.Count

# Shared-string insertion order observed in the target file is column-major
# (all of B, then all of C, then all of D) - A/E/F/G reuse numbers or
# already-existing strings, so write column by column to reproduce it.
for ($i = 0; $i -lt $n; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $A[$i]
}
for ($i = 0; $i -lt $n; $i++) {
    $ws.Cells.Item($i + 2, 2).Value = $B[$i]
}
for ($i = 0; $i -lt $n; $i++) {
    $ws.Cells.Item($i + 2, 3).Value = $C[$i]
}
for ($i = 0; $i -lt $n; $i++) {
    $ws.Cells.Item($i + 2, 4).Value = $D[$i]
}
for ($i = 0; $i -lt $n; $i++) {
    $ws.Cells.Item($i + 2, 5).Value = $E[$i]
}
for ($i = 0; $i -lt $n; $i++) {
    $ws.Cells.Item($i + 2, 6).Value = $F[$i]
}
for ($i = 0; $i -lt $n; $i++) {
    $ws.Cells.Item($i + 2, 7).Value = $G[$i]
}
